$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("B2").Value = 0.2093023255813954
$ws.Range("C2").Value = 0.5348837209302325
$ws.Range("P2").Value = 0.1395348837209302
$ws.Range("S2").Value = 0.1162790697674419
$ws.Range("J3").Value = 0.04166666666666666
$ws.Range("P3").Value = 0.75
$ws.Range("S3").Value = 0.2083333333333333
$ws.Range("P5").Value = 1
$ws.Range("B6").Value = 0.1578947368421053
$ws.Range("J6").Value = 0.3684210526315789
$ws.Range("O6").Value = 0.05263157894736842
$ws.Range("Q6").Value = 0.1578947368421053
$ws.Range("R6").Value = 0.05263157894736842
$ws.Range("S6").Value = 0.2105263157894737
$ws.Range("B7").Value = 0.125
$ws.Range("Q7").Value = 0.625
$ws.Range("S7").Value = 0.25
$ws.Range("B8").Value = 0.09375
$ws.Range("F8").Value = 0.03125
$ws.Range("J8").Value = 0.21875
$ws.Range("Q8").Value = 0.1875
$ws.Range("R8").Value = 0.03125
$ws.Range("S8").Value = 0.4375
$ws.Range("Q9").Value = 0.2
$ws.Range("R9").Value = 0.2
$ws.Range("B10").Value = 0.1982758620689655
$ws.Range("D10").Value = 0.01724137931034483
$ws.Range("E10").Value = 0.008620689655172414
$ws.Range("F10").Value = 0.0603448275862069
$ws.Range("J10").Value = 0.1379310344827586
$ws.Range("O10").Value = 0.01724137931034483
$ws.Range("Q10").Value = 0.1982758620689655
$ws.Range("R10").Value = 0.06896551724137931
$ws.Range("S10").Value = 0.293103448275862
$ws.Range("G11").Value = 0.1
$ws.Range("J11").Value = 0.1
$ws.Range("K11").Value = 0.1
$ws.Range("L11").Value = 0.6
$ws.Range("S11").Value = 0.1
$ws.Range("G12").Value = 0.75
$ws.Range("L12").Value = 0.25
$ws.Range("G13").Value = 0.3333333333333333
$ws.Range("J13").Value = 0.6666666666666666
$ws.Range("J15").Value = 0.4666666666666667
$ws.Range("K15").Value = 0.06666666666666667
$ws.Range("M15").Value = 0.06666666666666667
$ws.Range("S15").Value = 0.4
$ws.Range("H16").Value = 0.2592592592592592
$ws.Range("I16").Value = 0.07407407407407407
$ws.Range("J16").Value = 0.3703703703703703
$ws.Range("K16").Value = 0.03703703703703703
$ws.Range("O16").Value = 0.1111111111111111
$ws.Range("S16").Value = 0.1481481481481481
$ws.Range("F17").Value = 0.05
$ws.Range("H17").Value = 0.175
$ws.Range("I17").Value = 0.025
$ws.Range("J17").Value = 0.475
$ws.Range("K17").Value = 0.025
$ws.Range("O17").Value = 0.1
$ws.Range("S17").Value = 0.15
$ws.Range("H18").Value = 0.2727272727272727
$ws.Range("I18").Value = 0.1818181818181818
$ws.Range("J18").Value = 0.3636363636363636
$ws.Range("M18").Value = 0.09090909090909091
$ws.Range("S18").Value = 0.09090909090909091
$ws.Range("F19").Value = 0.05319148936170213
$ws.Range("H19").Value = 0.1595744680851064
$ws.Range("I19").Value = 0.05319148936170213
$ws.Range("J19").Value = 0.4680851063829787
$ws.Range("K19").Value = 0.06382978723404255
$ws.Range("M19").Value = 0.01063829787234043
$ws.Range("O19").Value = 0.05319148936170213
$ws.Range("S19").Value = 0.1382978723404255
